# Remove trailing spaces from header cell text values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "AI Analysis"
$ws.Range("F1").Value = "Planned Transition Partner"

# Update the active selection to I1 (matches the saved sheet view state)
$ws.Range("I1").Select()
